# Lunchtimes and lecturer time limits implemented.
# Replace the volatile RANDBETWEEN(0,1) occupancy formulas in B3:H67 with a
# flat literal "occupied" (1) baseline, and update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Collapse every formula cell in B3:H67 down to a plain literal value of 1
# (this both clears the shared RANDBETWEEN formulas and sets every value to 1).
$ws.Range("B3:H67").Value = 1

# Match the saved selection state from the edit.
$ws.Range("D74").Select()
